$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2022" year column is being added to the right of the existing
# "2021" column (S). Insert it so it inherits column S's formatting
# (number format / font / borders), just like the rest of the year columns.
$ws.Columns("T").Insert(-4121)

# Fill in the 2022 figures.
$ws.Range("T4").Value = 2022
$ws.Range("T5").Value = 3.7
$ws.Range("T6").Value = 1.6
$ws.Range("T7").Value = 1.7
$ws.Range("T8").Value = 17.899999999999999
$ws.Range("T9").Value = 7.5
$ws.Range("T10").Value = 1.1000000000000001
$ws.Range("T11").Value = 4.4000000000000004
$ws.Range("T12").Value = 3
$ws.Range("T13").Value = 4.0999999999999996
$ws.Range("T14").Value = 0.8

# Match the saved UI selection state.
$ws.Range("U4").Select()
